$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 header: "Empleado" -> "Legajo"
$ws.Range("A1").Value = "Legajo"

# Two new trailing columns: O "Motivo", P "Observación"
$ws.Cells.Item(1, 15).Value = "Motivo"
$ws.Cells.Item(1, 16).Value = "Observación"

# Give the new header cells the same look (bold/centered/bordered) as the
# rest of row 1 by copying the format from the last existing header cell.
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Materialize O2:P5 as (blank) cells under the two new columns for every
# existing data row, matching the new columns added to the table.
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 15).Borders.LineStyle = 0
    $ws.Cells.Item($r, 16).Borders.LineStyle = 0
}
